{"js": "// Edit: a second space was typed after \"digitalizing \" in the introduction\n// paragraph, so \"In a rapidly digitalizing world, ...\" became\n// \"In a rapidly digitalizing  world, ...\" (note the doubled space).\n//\n// We locate the unique \"digitalizing \" occurrence in the document body and\n// insert one extra space immediately after it (i.e. right before \"world\"),\n// mirroring a user placing the cursor there and pressing the space bar once.\n\nconst body = context.document.body;\n\nconst results = body.search(\"digitalizing \", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Could not find \"digitalizing \" in the document body.');\n}\n\n// Collapse the found range to its end (just before \"world\") and insert a\n// single space there, turning the single space into a double space.\nconst insertionPoint = results.items[0].getRange(Word.RangeLocation.end);\ninsertionPoint.insertText(\" \", Word.InsertLocation.before);\n\nawait context.sync();\n", "ps1": "# Edit: a second space was typed after \"digitalizing \" in the introduction\n# paragraph, so \"In a rapidly digitalizing world, ...\" became\n# \"In a rapidly digitalizing  world, ...\" (note the doubled space).\n#\n# We find the unique \"digitalizing \" occurrence in the document and insert a\n# single extra space right after it (i.e. immediately before \"world\"),\n# mirroring a user placing the cursor there and pressing the space bar once.\n\n$d = $word.ActiveDocument\n\n$r = $d.Content\n$find = $r.Find\n$find.ClearFormatting()\n$find.Text = \"digitalizing \"\n$find.MatchCase = $true\n$found = $find.Execute()\n\nif (-not $found) {\n    throw 'Could not find \"digitalizing \" in the document.'\n}\n\n# $r now spans the matched text (\"digitalizing \"); collapse it to its end\n# point (just before \"world\") and insert one space there.\n$r.Collapse(0)  # wdCollapseEnd\n$r.InsertAfter(\" \")\n"}
